$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Shartnoma raqam) and G (Telefon raqam) to be stored as text
# for the new rows, matching the existing data convention in the sheet.
$ws.Range("D73:D87").NumberFormat = "@"
$ws.Range("G73:G87").NumberFormat = "@"

$ws.Range("A73").Value2 = 'Rasulova Odina Abduvaliyevna'
$ws.Range("B73").Value2 = 'Amaliy psixologiya 576 soatlik'
$ws.Range("C73").Value2 = 'AD3305485'
$ws.Range("D73").Value2 = '669'
$ws.Range("E73").Value2 = 'Fargona viloyati'
$ws.Range("F73").Value2 = 'Rishton tumani'
$ws.Range("G73").Value2 = '998905633728'
$ws.Range("H73").Value2 = '25-10-2024'

$ws.Range("A74").Value2 = 'Alimova Hayotxon Qobiljonovna'
$ws.Range("B74").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C74").Value2 = 'AD5159476'
$ws.Range("D74").Value2 = '670'
$ws.Range("E74").Value2 = 'Andijon viloyati'
$ws.Range("F74").Value2 = 'Jalaquduq tuman'
$ws.Range("G74").Value2 = '998916126213'
$ws.Range("H74").Value2 = '25-10-2024'

$ws.Range("A75").Value2 = 'Madg''oziyeva Sabohat Xatamovna'
$ws.Range("B75").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C75").Value2 = 'AD3227666'
$ws.Range("D75").Value2 = '671'
$ws.Range("E75").Value2 = 'Fargona viloyati'
$ws.Range("F75").Value2 = 'Toshloq tumani'
$ws.Range("G75").Value2 = '998916580488'
$ws.Range("H75").Value2 = '25-10-2024'

$ws.Range("A76").Value2 = 'Omonkeldiyeva Nilufar Shuxrat qizi'
$ws.Range("B76").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C76").Value2 = 'AD1581888'
$ws.Range("D76").Value2 = '672'
$ws.Range("E76").Value2 = 'Andijon viloyati'
$ws.Range("F76").Value2 = 'Andijon tuman'
$ws.Range("G76").Value2 = '998915227679'
$ws.Range("H76").Value2 = '25-10-2024'

$ws.Range("A77").Value2 = 'Xaydarova Dilfuza Ikromjon qizi'
$ws.Range("B77").Value2 = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C77").Value2 = 'AD0898964'
$ws.Range("D77").Value2 = '673'
$ws.Range("E77").Value2 = 'Fargona viloyati'
$ws.Range("F77").Value2 = 'Oltiariq tumani'
$ws.Range("G77").Value2 = '998916642168'
$ws.Range("H77").Value2 = '26-10-2024'

$ws.Range("A78").Value2 = 'Subhonova Farida Islomiddin qizi'
$ws.Range("B78").Value2 = 'Amaliy psixologiya 576 soatlik'
$ws.Range("C78").Value2 = 'AB4455226'
$ws.Range("D78").Value2 = '674'
$ws.Range("E78").Value2 = 'Samarqand viloyati'
$ws.Range("F78").Value2 = 'Pastdargʻom tumani'
$ws.Range("G78").Value2 = '998946838894'
$ws.Range("H78").Value2 = '26-10-2024'

$ws.Range("A79").Value2 = 'Sadullayeva Nodira Saparbay qizi'
$ws.Range("B79").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C79").Value2 = 'AB9081455'
$ws.Range("D79").Value2 = '675'
$ws.Range("E79").Value2 = 'Xorazm viloyati'
$ws.Range("F79").Value2 = 'Yangibozor tumani'
$ws.Range("G79").Value2 = '998943296838'
$ws.Range("H79").Value2 = '26-10-2024'

$ws.Range("A80").Value2 = 'Elboyeva Nazokat Otabek qizi'
$ws.Range("B80").Value2 = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C80").Value2 = 'AB7009475'
$ws.Range("D80").Value2 = '676'
$ws.Range("E80").Value2 = 'Navoiy viloyati'
$ws.Range("F80").Value2 = 'Karmana tumani'
$ws.Range("G80").Value2 = '998935614006'
$ws.Range("H80").Value2 = '28-10-2024'

$ws.Range("A81").Value2 = 'Yuldosheva Oybibi Selxanovna'
$ws.Range("B81").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C81").Value2 = 'AB9839541'
$ws.Range("D81").Value2 = '677'
$ws.Range("E81").Value2 = 'Xorazm viloyati'
$ws.Range("F81").Value2 = 'Urganch tumani'
$ws.Range("G81").Value2 = '998993808528'
$ws.Range("H81").Value2 = '28-10-2024'

$ws.Range("A82").Value2 = 'Xolillayeva Kumushoy Ikromjon qizi'
$ws.Range("B82").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C82").Value2 = 'AD4560577'
$ws.Range("D82").Value2 = '678'
$ws.Range("E82").Value2 = 'Qoraqalpogʻiston Respublikasi'
$ws.Range("F82").Value2 = 'Amudaryo tumani'
$ws.Range("G82").Value2 = '998971705770'
$ws.Range("H82").Value2 = '28-10-2024'

$ws.Range("A83").Value2 = 'Normamatov Hamrozbek Xusniddin o''g''li'
$ws.Range("B83").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C83").Value2 = 'AC0380497'
$ws.Range("D83").Value2 = '679'
$ws.Range("E83").Value2 = 'Samarqand viloyati'
$ws.Range("F83").Value2 = 'Kattaqoʻrgʻon tumani'
$ws.Range("G83").Value2 = '998936066607'
$ws.Range("H83").Value2 = '29-10-2024'

$ws.Range("A84").Value2 = 'Kiyasova Ulmeken Ilyas qizi'
$ws.Range("B84").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C84").Value2 = 'KA0680487'
$ws.Range("D84").Value2 = '680'
$ws.Range("E84").Value2 = 'Qoraqalpogʻiston Respublikasi'
$ws.Range("F84").Value2 = 'Qoʻngʻirot tumani'
$ws.Range("G84").Value2 = '998994554845'
$ws.Range("H84").Value2 = '29-10-2024'

$ws.Range("A85").Value2 = 'Axunova Shohzodaxon Olim qizi'
$ws.Range("B85").Value2 = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C85").Value2 = 'AD1635577'
$ws.Range("D85").Value2 = '681'
$ws.Range("E85").Value2 = 'Toshkent shahri'
$ws.Range("F85").Value2 = 'Yunusobod tumani'
$ws.Range("G85").Value2 = '998998009889'
$ws.Range("H85").Value2 = '30-10-2024'

$ws.Range("A86").Value2 = 'Haydarova Maftuna Haliljon qizi'
$ws.Range("B86").Value2 = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C86").Value2 = 'AD6997411'
$ws.Range("D86").Value2 = '682'
$ws.Range("E86").Value2 = 'Fargona viloyati'
$ws.Range("F86").Value2 = 'Oltiariq tumani'
$ws.Range("G86").Value2 = '998931690141'
$ws.Range("H86").Value2 = '30-10-2024'

$ws.Range("A87").Value2 = 'Muhtorova Nargiza Abduvositovna'
$ws.Range("B87").Value2 = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C87").Value2 = 'AD5529844'
$ws.Range("D87").Value2 = '683'
$ws.Range("E87").Value2 = 'Andijon viloyati'
$ws.Range("F87").Value2 = 'Andijon tuman'
$ws.Range("G87").Value2 = '998999013032'
$ws.Range("H87").Value2 = '30-10-2024'
